$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: mirrors row 5 pattern (TC1 login, no extra detail columns)
$ws.Range("A8").Value = "infor_TC1_Login2"
$ws.Range("B8").Value = "infor_taas"
$ws.Range("C8").Value = "infor_taas"

# Row 9: mirrors row 6 pattern (TC2 registration, with full detail columns + hyperlink)
$ws.Range("A9").Value = "infor_TC2_Registration2"
$ws.Range("B9").Value = "infor_taas"
$ws.Range("C9").Value = "infor_taas"
$ws.Range("D9").Value = "Zenric"
$ws.Range("E9").Value = "Navea"
$ws.Range("F9").Value = "jazx.zn@gmail.com"
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:jazx.zn@gmail.com")
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("G9").Value = "'09661401029"
$ws.Range("H9").Value = "Male"
$ws.Range("I9").Value = "RCG Technology Inc."
$ws.Range("J9").Value = "Address Line 1"
$ws.Range("K9").Value = "Address Line 2"
$ws.Range("L9").Value = "Automation"
$ws.Range("M9").Value = "Selenium WebDriver"

# Row 10: mirrors row 7 pattern (TC3 login via config file, no extra detail columns)
$ws.Range("A10").Value = "infor_TC3_LoginViaConfigFile2"
$ws.Range("B10").Value = "infor_taas"
$ws.Range("C10").Value = "infor_taas"

# Move the active selection to A11, matching the post-edit cursor position
$ws.Range("A11").Select()
